$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the erroring INDEX formula in D6 with a working one
$ws.Range("D6").Formula = "=INDEX(A1:B5,1,2)"

# Add two new simple arithmetic formulas in D7 and D8
$ws.Range("D7").Formula = "=2+3"
$ws.Range("D8").Formula = "=3+2"

$wb.Save()
